$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample data rows (rows 2-4); keep the header row (row 1)
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Add the new header columns F1:I1
$ws.Range("F1").Value = "Warnings"
$ws.Range("G1").Value = "Home Address"
$ws.Range("H1").Value = "Balance"
$ws.Range("I1").Value = "Phone number"

# Copy the existing header formatting onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize the columns to the new layout widths
$ws.Columns.Item(2).ColumnWidth = 15.709635416666666
$ws.Columns.Item(3).ColumnWidth = 18.436197916666668
$ws.Columns.Item(4).ColumnWidth = 10.529947916666666
$ws.Columns.Item(5).ColumnWidth = 17.529947916666668
$ws.Columns.Item(7).ColumnWidth = 14.072916666666666
$ws.Columns.Item(8).ColumnWidth = 10.529947916666666
$ws.Columns.Item(9).ColumnWidth = 13.436197916666666

# Move the active selection to A2
[void]$ws.Range("A2").Select()
